$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.338.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.53%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.932.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.67%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'250.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.85%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.7162"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.62%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.25%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3272"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.97%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'27.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.97%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.8011"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.12%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08075"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.91%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.933.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.75%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.420"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.80%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'94.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.29%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +3.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.334.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.47%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'252.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.56%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000008127"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.13%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'5.810"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.186.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.90%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.30%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.56%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.934"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'9.724"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.97%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'165.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.39%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.334"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.99%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'19.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.74%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.1291"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.73%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.68%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.541"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.24%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.422"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.198"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.51%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.05203"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.65%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.266"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7469"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.09%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.763"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.15%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.16%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'79.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.97%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.464"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.17%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4520"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.63%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.024"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.82%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.20%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.8393"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.16%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'101.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.32%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'9.791"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.87%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.412"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.46%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'36.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.92%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06062"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.71%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.4174"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.29%  "
$ws.Range("E51").Style = "Normal"
